$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing hyperlinks will be stale once rows shift down, so drop them now
# and re-add them (in column order) once every row is in its final place.
$ws.Hyperlinks.Delete()

# Insert a new row at position 5 - this pushes the old rows 5-7 down to 6-8
# and carries their cell values/formatting along with them.
$ws.Rows(5).Insert()

# Fill in the brand-new row 5 with the newly scraped job posting.
$ws.Range("A5").Value = "2025-11-24 01:57:09"
$ws.Range("B5").Value = "マッチングサイト開発エンジニア募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5440077"
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = "◆開発 ◇サイト"

# This scrape refreshed every row's "取得日時" (fetched-at) timestamp.
$ws.Range("A2").Value = "2025-11-24 01:57:09"
$ws.Range("A3").Value = "2025-11-24 01:57:09"
$ws.Range("A4").Value = "2025-11-24 01:57:09"
$ws.Range("A6").Value = "2025-11-24 01:57:09"
$ws.Range("A7").Value = "2025-11-24 01:57:09"
$ws.Range("A8").Value = "2025-11-24 01:57:09"

# Re-create the URL hyperlinks now that every row sits in its final slot.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5440052")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5439921")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5440010")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5440077")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5439670")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5440042")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5440002")

# Hyperlinks.Add() stamps a brand-new style; reapply the workbook's existing
# "Hyperlink" cell style so the cells match the original formatting exactly.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
